$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update ConceptScheme URI (B1) and PREFIX URI (C3)
$ws.Range("B1").Value = "http://purl.org/m4m21/variables"
$ws.Range("C3").Value = "http://purl.org/m4m21/variables/"

# Update last-modified timestamp
$ws.Range("B20").Value = "2022-06-05T08:36:30+00:00"

# Row 23: rename top term, clear alt label / definition
$ws.Range("B23").Value = "new top variable"
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = ""

# Row 24: rename narrow term and update its broader reference
$ws.Range("B24").Value = "new narrow variable"
$ws.Range("F24").Value = "new top variable"

# Remove rows 25 and 26 (Test Top Subject 2 / Test Narrow Subject 2) entirely
$ws.Range("A25:AK26").EntireRow.Delete()
